# Update the "想去人数" (F) and "最低票价" (G) figures across the
# "展览" (sheet1), "本地生活" (sheet3) and "全部类型" (sheet4) worksheets
# to reflect the newly generated gh-pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 21432
$ws1.Range("F3").Value  = 3362
$ws1.Range("G3").Value  = 70
$ws1.Range("G4").Value  = 70
$ws1.Range("F9").Value  = 268
$ws1.Range("F10").Value = 75
$ws1.Range("F12").Value = 588
$ws1.Range("F14").Value = 369
$ws1.Range("F15").Value = 42
$ws1.Range("F16").Value = 469
$ws1.Range("F17").Value = 227
$ws1.Range("F18").Value = 46
$ws1.Range("F21").Value = 162

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6182
$ws3.Range("F5").Value = 1748
$ws3.Range("F6").Value = 95

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6182
$ws4.Range("F5").Value  = 1748
$ws4.Range("F6").Value  = 21432
$ws4.Range("F7").Value  = 3362
$ws4.Range("G7").Value  = 70
$ws4.Range("G8").Value  = 70
$ws4.Range("F10").Value = 95
$ws4.Range("F15").Value = 268
$ws4.Range("F17").Value = 75
$ws4.Range("F23").Value = 588
$ws4.Range("F27").Value = 369
$ws4.Range("F29").Value = 42
$ws4.Range("F30").Value = 469
$ws4.Range("F32").Value = 227
$ws4.Range("F33").Value = 46
$ws4.Range("F43").Value = 162
